# Ajout du nom de projet
#
# Fills in the project name on the "Groupe" sheet, next to the
# "Nom du produit :" label in B27 (mirrors the "Nom :" / student-name
# pattern used in column C of the "Etudiant N" sheets).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Groupe")

$ws.Range("C27").Value = "ProjetFramboise"

# Keep the recorded selection in step with the edit (matches the
# author's final cursor position after typing the project name).
$ws.Range("AB35").Select() | Out-Null
